# Updated RAD test cases for Existing Liability - Motor Fuel Tax
# The "Date" column (B) logs the timestamp at which each RAD test step
# was recorded. This run refreshes those timestamps for rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sun Oct 13 00:05:35 EDT 2024"
$ws.Range("B3").Value = "Sun Oct 13 00:05:47 EDT 2024"
$ws.Range("B4").Value = "Sun Oct 13 00:06:00 EDT 2024"
$ws.Range("B5").Value = "Sun Oct 13 00:06:12 EDT 2024"
$ws.Range("B6").Value = "Sun Oct 13 00:06:24 EDT 2024"
$ws.Range("B7").Value = "Sun Oct 13 00:06:37 EDT 2024"
